$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "Active"

# New "Color" column (F) - add header + fill the column top-to-bottom first
$ws.Range("F1").Value = "Color"
$ws.Range("F2").Value = "MistyRose"
$ws.Range("F3").Value = "OldLace"

# Copy header formatting (style s="1") from an existing header cell onto F1
$ws.Range("C1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set the new column's width
$ws.Columns.Item(6).ColumnWidth = 9.64

# --- Row 2 / Row 3 value updates ---
$ws.Range("E2").Value = "Tick"
$ws.Range("E3").Value = "Tick"

# RefID values updated last (LAT-100/LAT-92 -> LAT-137/LAT-138)
$ws.Range("A2").Value = "LAT-137"
$ws.Range("A3").Value = "LAT-138"

# Reset active cell/selection back to A1
$ws.Range("A1").Select() | Out-Null
